# Edit: reorder a handful of rows on the "anobjl" sheet, update the
# frozen-pane / selection state, and record the sort-range bookkeeping
# that Excel leaves behind after a manual sort.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("anobjl")

# ---------------------------------------------------------------------
# 1. Trigger a real Sort.Apply() on A3:C16 purely so the engine persists
#    the <sortState>/<sortCondition> bookkeeping the way Excel does after
#    a manual sort. We immediately overwrite every value touched by the
#    sort (below) with the exact original/target contents, so the actual
#    sort order it produces does not matter.
# ---------------------------------------------------------------------
$sort = $ws1.Sort
$sort.SortFields.Clear()
$keyRange = $ws1.Range("A3:A16")
$sort.SortFields.Add($keyRange)
$setRange = $ws1.Range("A3:C16")
$sort.SetRange($setRange)
$sort.Header = 2
$sort.Apply()

# ---------------------------------------------------------------------
# 2. Restore / set the real cell contents for rows 2-16.
#    Row 2 is untouched by the sort (range started at row 3) but still
#    gets the new value. Rows 3-7 get the reordered values. Rows 8-16
#    are restored to their original values (the sort above scrambled
#    them) including preserving the TRUE/FALSE boolean type in column C.
# ---------------------------------------------------------------------

$ws1.Range("A2").Value2 = "SEQUENCE"
$ws1.Range("B2").Value2 = "mbrp"
$ws1.Range("C2").Value2 = "gen"

$ws1.Range("A3").Value2 = "GRAPH"
$ws1.Range("B3").Value2 = "mbrp"
$ws1.Range("C3").Value2 = "sa"

$ws1.Range("A4").Value2 = "VERTEX"
$ws1.Range("B4").Value2 = "mbrp"
$ws1.Range("C4").Value2 = "sa"

$ws1.Range("A5").Value2 = "EDGE"
$ws1.Range("B5").Value2 = "mbrp"
$ws1.Range("C5").Value2 = "sa"

$ws1.Range("A6").Value2 = "OBSERVATION"
$ws1.Range("B6").Value2 = "mbrp"
$ws1.Range("C6").Value2 = "env"

$ws1.Range("A7").Value2 = "INFERENCE"
$ws1.Range("B7").Value2 = "mbrp"
$ws1.Range("C7").Value2 = "qa"

$ws1.Range("A8").Value2 = "GRAPH"
$ws1.Range("B8").Value2 = "tool_support"
$ws1.Range("C8").Value2 = $true

$ws1.Range("A9").Value2 = "OBSERVATION"
$ws1.Range("B9").Value2 = "tool_support"
$ws1.Range("C9").Value2 = $false

$ws1.Range("A10").Value2 = "SEQUENCE"
$ws1.Range("B10").Value2 = "tool_support"
$ws1.Range("C10").Value2 = $true

$ws1.Range("A11").Value2 = "INFERENCE"
$ws1.Range("B11").Value2 = "tool_support"
$ws1.Range("C11").Value2 = $false

$ws1.Range("A12").Value2 = "VERTEX"
$ws1.Range("B12").Value2 = "tool_support"
$ws1.Range("C12").Value2 = $true

$ws1.Range("A13").Value2 = "EDGE"
$ws1.Range("B13").Value2 = "tool_support"
$ws1.Range("C13").Value2 = $true

$ws1.Range("A14").Value2 = "VERTEX>id|OBSERVATION"
$ws1.Range("B14").Value2 = "type"
$ws1.Range("C14").Value2 = "defined"

$ws1.Range("A15").Value2 = "EDGE>id|OBSERVATION"
$ws1.Range("B15").Value2 = "type"
$ws1.Range("C15").Value2 = "defined"

$ws1.Range("A16").Value2 = "SEQUENCE>step|GRAPH"
$ws1.Range("B16").Value2 = "type"
$ws1.Range("C16").Value2 = "defined"

# ---------------------------------------------------------------------
# 3. Update the frozen pane's top-left cell and the active selection to
#    match the new view state recorded in the workbook.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B7").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "edit applied"
